$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'57.166.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  -1.06%  "
$ws.Range("D3").Value2 = "'2.988.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -2.22%  "
$ws.Range("E4").Value2 = "  +0.12%  "
$ws.Range("D5").Value2 = "'500.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -4.55%  "
$ws.Range("D6").Value2 = "'138.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -3.01%  "
$ws.Range("E7").Value2 = "  +0.00%  "
$ws.Range("E8").Value2 = "  -3.48%  "
$ws.Range("D9").Value2 = "'7.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  -4.25%  "
$ws.Range("E10").Value2 = "  -4.65%  "
$ws.Range("E11").Value2 = "  -3.17%  "
$ws.Range("D12").Value2 = "'3.505.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -2.02%  "
$ws.Range("E13").Value2 = "  -2.30%  "
$ws.Range("D14").Value2 = "'26.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -2.62%  "
$ws.Range("E15").Value2 = "  -6.38%  "
$ws.Range("D16").Value2 = "'57.220.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -0.93%  "
$ws.Range("E17").Value2 = "  -2.54%  "
$ws.Range("D18").Value2 = "'2.989.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -2.07%  "
$ws.Range("D19").Value2 = "'12.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -3.18%  "
$ws.Range("D20").Value2 = "'7.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -3.73%  "
$ws.Range("D21").Value2 = "'320.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -5.50%  "
$ws.Range("D22").Value2 = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -0.14%  "
$ws.Range("D23").Value2 = "'5.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +0.77%  "
$ws.Range("E24").Value2 = "  -1.98%  "
$ws.Range("D25").Value2 = "'63.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -2.29%  "
$ws.Range("E26").Value2 = "  -0.25%  "
$ws.Range("E27").Value2 = "  -5.62%  "
$ws.Range("D28").Value2 = "'0.0₃0895"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -8.26%  "
$ws.Range("E29").Value2 = "  -5.01%  "
$ws.Range("E30").Value2 = "  -3.27%  "
$ws.Range("E31").Value2 = "  -4.10%  "
$ws.Range("E32").Value2 = "  -4.92%  "
$ws.Range("D33").Value2 = "'20.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -4.23%  "
$ws.Range("D34").Value2 = "'155.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -0.61%  "
$ws.Range("D35").Value2 = "'4.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -3.92%  "
$ws.Range("D36").Value2 = "'5.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -3.27%  "
$ws.Range("D37").Value2 = "'1.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -5.82%  "
$ws.Range("D38").Value2 = "'24.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -6.47%  "
$ws.Range("D39").Value2 = "'0.0663"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -5.83%  "
$ws.Range("D40").Value2 = "'3.023.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -2.25%  "
$ws.Range("D41").Value2 = "'37.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +0.26%  "
$ws.Range("D42").Value2 = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -0.07%  "
$ws.Range("E43").Value2 = "  -3.65%  "
$ws.Range("E44").Value2 = "  -2.68%  "
$ws.Range("E45").Value2 = "  -6.13%  "
$ws.Range("D46").Value2 = "'2.196.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -5.76%  "
$ws.Range("D47").Value2 = "'5.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -1.85%  "
$ws.Range("D48").Value2 = "'0.939"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -9.27%  "
$ws.Range("E49").Value2 = "  -4.83%  "
$ws.Range("D50").Value2 = "'19.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -4.50%  "
$ws.Range("D51").Value2 = "'1.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -11.92%  "
